$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: convert SmartScore text cells to real numbers ---
$ws.Range("G4").Value = 0.572
$ws.Range("J4").Value = 0.488
$ws.Range("M4").Value = 0.412
$ws.Range("P4").Value = 0.636
$ws.Range("S4").Value = 0.578
$ws.Range("V4").Value = 0.541
$ws.Range("Y4").Value = 0.738
$ws.Range("AB4").Value = 0.614
$ws.Range("AE4").Value = 0.599

# --- Row 5: new submission (Miranda) ---
$ws.Range("A5").Value = 'Miranda'
$ws.Range("B5").Value = 25
$ws.Range("C5").Value = 'Femenino'
$ws.Range("D5").Value = '2025-10-28 05:51:45'
$ws.Range("E5").Value = @'
{
  "portion": 0.8,
  "diet": 0.5714285714285714,
  "salt": 0.6,
  "fat": 0.8,
  "natural": 0.6,
  "convenience": 0.4,
  "price": 0.8
}
'@
$ws.Range("F5").Value = 'Nongshim Neoguri Spicy Seafood'
$ws.Range("G5").Value = '''0.575'
$ws.Range("H5").Value = 'Sabor a marisco, umami, picante equilibrado, buena textura, algo salado'
$ws.Range("I5").Value = 'Nissin Chow Mein Teriyaki Beef'
$ws.Range("J5").Value = '''0.510'
$ws.Range("K5").Value = 'Fácil de preparar, porción generosa, salsa suave, necesita mejoras, alto en grasa'
$ws.Range("L5").Value = 'Maruchan Ramen Sabor Pollo'
$ws.Range("M5").Value = '''0.509'
$ws.Range("N5").Value = 'Sabor clásico, económico, alto en sodio, no saludable, nostálgico'
$ws.Range("O5").Value = 'Kraft Macaroni & Cheese Dinner'
$ws.Range("P5").Value = '''0.650'
$ws.Range("Q5").Value = 'Sabor nostálgico, clásico americano, fácil, no muy nutritivo, barato'
$ws.Range("R5").Value = 'Annie’s Shells & White Cheddar'
$ws.Range("S5").Value = '''0.587'
$ws.Range("T5").Value = 'Queso blanco real, sin colorantes, sabor casero, menos salado, buena para niños'
$ws.Range("U5").Value = 'Amy’s Macaroni & Cheese (frozen)'
$ws.Range("V5").Value = '''0.552'
$ws.Range("W5").Value = 'Queso real, textura casera, sin conservadores, alto en grasa, algo caro'
$ws.Range("X5").Value = 'Wild Planet Wild Tuna Pasta Salad'
$ws.Range("Y5").Value = '''0.664'
$ws.Range("Z5").Value = 'Sabor fresco, buena proteína, saludable, porción algo pequeña'
$ws.Range("AA5").Value = 'StarKist Chicken Creations (Chicken Salad)'
$ws.Range("AB5").Value = '''0.589'
$ws.Range("AC5").Value = 'Portátil, saludable, fácil, buena textura, sabor suave'
$ws.Range("AD5").Value = 'Jack Link’s Beef Jerky Original'
$ws.Range("AE5").Value = '''0.576'
$ws.Range("AF5").Value = 'Ahumado, sabroso, alto en proteína, snack ideal, porción pequeña'

# Reset row height (the multi-line JSON in E5 otherwise triggers an auto row-height
# bump that the source workbook does not have).
$ws.Rows.Item(5).AutoFit()
